$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage-formatted text cells (column H) need an explicit Text number format
# before assignment, otherwise Excel auto-converts "NN%" strings into numeric
# percentage values instead of leaving them as literal text.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "46%"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "88%"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "48%"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "84%"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "42%"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "42%"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "62%"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "50%"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "51%"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "51%"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "71%"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "39%"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "30%"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "42%"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "92%"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "50%"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "43%"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "56%"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "89%"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "60%"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "51%"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "51%"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "53%"

# Remaining text cells (dates, pressures, temperatures, wind, radiation, etc.)
$ws.Range("E2").Value = "2026-02-17 07:18:47"
$ws.Range("E3").Value = "2026-02-17 07:18:50"
$ws.Range("O3").Value = "-6.9 °C"
$ws.Range("E4").Value = "2026-02-17 07:18:52"
$ws.Range("J4").Value = "1016.0 hPa"
$ws.Range("M4").Value = "6.9 °C 6:46 TU"
$ws.Range("E5").Value = "2026-02-17 07:18:55"
$ws.Range("O5").Value = "-6.3 °C"
$ws.Range("E6").Value = "2026-02-17 07:18:57"
$ws.Range("J6").Value = "1015.7 hPa"
$ws.Range("E7").Value = "2026-02-17 07:19:00"
$ws.Range("J7").Value = "1015.3 hPa"
$ws.Range("N7").Value = "12.4 °C 6:59 TU"
$ws.Range("O7").Value = "13.8 °C"
$ws.Range("E8").Value = "2026-02-17 07:19:02"
$ws.Range("J8").Value = "1015.5 hPa"
$ws.Range("N8").Value = "8.5 °C 6:59 TU"
$ws.Range("E9").Value = "2026-02-17 07:19:04"
$ws.Range("O9").Value = "12.1 °C"
$ws.Range("E10").Value = "2026-02-17 07:19:07"
$ws.Range("K10").Value = "0.0 MJ/m2"
$ws.Range("L10").Value = "18.4 km/h - 103º 6:45 TU"
$ws.Range("O10").Value = "8.2 °C"
$ws.Range("E11").Value = "2026-02-17 07:19:09"
$ws.Range("N11").Value = "1.8 °C 6:34 TU"
$ws.Range("O11").Value = "5.4 °C"
$ws.Range("E12").Value = "2026-02-17 07:19:12"
$ws.Range("O12").Value = "12.4 °C"
$ws.Range("E13").Value = "2026-02-17 07:19:14"
$ws.Range("J13").Value = "1017.2 hPa"
$ws.Range("O13").Value = "4.2 °C"
$ws.Range("E14").Value = "2026-02-17 07:19:15"
$ws.Range("N14").Value = "7.9 °C 6:51 TU"
$ws.Range("O14").Value = "12.4 °C"
$ws.Range("E15").Value = "2026-02-17 07:19:16"
$ws.Range("N15").Value = "9.8 °C 6:47 TU"
$ws.Range("O15").Value = "11.8 °C"
$ws.Range("E16").Value = "2026-02-17 07:19:18"
$ws.Range("E17").Value = "2026-02-17 07:19:21"
$ws.Range("K17").Value = "0.0 MJ/m2"
$ws.Range("N17").Value = "0.6 °C 6:55 TU"
$ws.Range("O17").Value = "2.5 °C"
$ws.Range("E18").Value = "2026-02-17 07:19:23"
$ws.Range("J18").Value = "1016.0 hPa"
$ws.Range("N18").Value = "4.5 °C 6:50 TU"
$ws.Range("O18").Value = "6.2 °C"
$ws.Range("E19").Value = "2026-02-17 07:19:25"
$ws.Range("N19").Value = "3.8 °C 6:53 TU"
$ws.Range("O19").Value = "5.8 °C"
$ws.Range("E20").Value = "2026-02-17 07:19:26"
$ws.Range("E21").Value = "2026-02-17 07:19:27"
$ws.Range("J21").Value = "1016.1 hPa"
$ws.Range("O21").Value = "7.5 °C"
$ws.Range("E22").Value = "2026-02-17 07:19:28"
$ws.Range("E23").Value = "2026-02-17 07:19:29"
$ws.Range("M23").Value = "-5.5 °C 6:59 TU"
$ws.Range("E24").Value = "2026-02-17 07:19:31"
$ws.Range("J24").Value = "1017.7 hPa"
$ws.Range("N24").Value = "9.1 °C 6:48 TU"
$ws.Range("E25").Value = "2026-02-17 07:19:32"
$ws.Range("E26").Value = "2026-02-17 07:19:34"
$ws.Range("E27").Value = "2026-02-17 07:19:36"
$ws.Range("E28").Value = "2026-02-17 07:19:39"
$ws.Range("J28").Value = "1016.2 hPa"
$ws.Range("O28").Value = "5.1 °C"
$ws.Range("E29").Value = "2026-02-17 07:19:41"
$ws.Range("K29").Value = "0.0 MJ/m2"
$ws.Range("O29").Value = "11.5 °C"
$ws.Range("E30").Value = "2026-02-17 07:19:44"
$ws.Range("J30").Value = "1015.3 hPa"
$ws.Range("K30").Value = "0.0 MJ/m2"
$ws.Range("N30").Value = "8.8 °C 6:31 TU"
$ws.Range("O30").Value = "11.5 °C"
$ws.Range("E31").Value = "2026-02-17 07:19:47"
$ws.Range("J31").Value = "1015.8 hPa"
$ws.Range("E32").Value = "2026-02-17 07:19:49"
$ws.Range("N32").Value = "5.4 °C 6:48 TU"
$ws.Range("E33").Value = "2026-02-17 07:19:51"
$ws.Range("J33").Value = "1016.2 hPa"
$ws.Range("O33").Value = "4.6 °C"
$ws.Range("E34").Value = "2026-02-17 07:19:54"
$ws.Range("M34").Value = "0.7 °C 6:42 TU"
$ws.Range("O34").Value = "-0.8 °C"
$ws.Range("E35").Value = "2026-02-17 07:19:56"
$ws.Range("J35").Value = "1018.6 hPa"
$ws.Range("O35").Value = "5.3 °C"
$ws.Range("E36").Value = "2026-02-17 07:19:58"
$ws.Range("J36").Value = "1015.8 hPa"
$ws.Range("O36").Value = "12.4 °C"
$ws.Range("E37").Value = "2026-02-17 07:20:01"
$ws.Range("J37").Value = "1016.3 hPa"
$ws.Range("N37").Value = "4.4 °C 6:35 TU"
$ws.Range("O37").Value = "7.2 °C"
$ws.Range("E38").Value = "2026-02-17 07:20:04"
$ws.Range("E39").Value = "2026-02-17 07:20:06"
$ws.Range("O39").Value = "-4.5 °C"
$ws.Range("E40").Value = "2026-02-17 07:20:09"
$ws.Range("J40").Value = "1017.8 hPa"
$ws.Range("O40").Value = "5.2 °C"
$ws.Range("E41").Value = "2026-02-17 07:20:11"
$ws.Range("J41").Value = "1015.9 hPa"
$ws.Range("O41").Value = "14.4 °C"
$ws.Range("E42").Value = "2026-02-17 07:20:13"
$ws.Range("O42").Value = "12.4 °C"
$ws.Range("E43").Value = "2026-02-17 07:20:16"
$ws.Range("K43").Value = "0.0 MJ/m2"
$ws.Range("O43").Value = "3.9 °C"
$ws.Range("E44").Value = "2026-02-17 07:20:18"
$ws.Range("E45").Value = "2026-02-17 07:20:20"
$ws.Range("J45").Value = "1021.1 hPa"
$ws.Range("E46").Value = "2026-02-17 07:20:23"
$ws.Range("J46").Value = "1018.0 hPa"
$ws.Range("N46").Value = "12.5 °C 6:59 TU"
$ws.Range("O46").Value = "13.3 °C"
